$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.609.02"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "1.754.96"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").Value = "  +0.41%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "323.95"
$ws.Range("E5").Value = "  +0.49%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  +0.41%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4509"
$ws.Range("E7").Value = "  +5.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3593"
$ws.Range("E8").Value = "  -0.52%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.96"
$ws.Range("E10").Value = "  -5.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.096"
$ws.Range("E11").Value = "  -1.44%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.9999"
$ws.Range("E12").Value = "  +0.43%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.88"
$ws.Range("E13").Value = "  -3.33%  "
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.108"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("D16").Value = "1.752.91"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "93.23"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("E19").Value = "  +0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9997"
$ws.Range("E20").Value = "  +0.40%  "
$ws.Range("E21").Value = "  -2.40%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.834"
$ws.Range("E22").Value = "  -2.89%  "
$ws.Range("D23").Value = "27.654.72"
$ws.Range("E23").Value = "  -1.41%  "
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.104"
$ws.Range("E25").Value = "  -2.35%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "162.66"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.47"
$ws.Range("E27").Value = "  +0.57%  "
$ws.Range("D28").Value = "1.954.13"
$ws.Range("E28").Value = "  -1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.081"
$ws.Range("E29").Value = "  -3.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.96"
$ws.Range("E30").Value = "  +0.47%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.085"
$ws.Range("E31").Value = "  -8.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09102"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.667"
$ws.Range("E33").Value = "  +4.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.553"
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "11.99"
$ws.Range("E35").Value = "  -5.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02294"
$ws.Range("E36").Value = "  -1.82%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2105"
$ws.Range("E37").Value = "  -0.50%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6387"
$ws.Range("E38").Value = "  -0.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06031"
$ws.Range("E39").Value = "  -1.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.960"
$ws.Range("E40").Value = "  -2.46%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.203"
$ws.Range("E41").Value = "  +1.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.384"
$ws.Range("E42").Value = "  -1.31%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.767"
$ws.Range("E43").Value = "  -1.52%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.33"
$ws.Range("E44").Value = "  -1.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5919"
$ws.Range("E45").Value = "  -1.30%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.712"
$ws.Range("E46").Value = "  +0.31%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "122.81"
$ws.Range("E47").Value = "  -0.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.955"
$ws.Range("E48").Value = "  -2.62%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.147"
$ws.Range("E49").Value = "  -2.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06862"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.34"
$ws.Range("E51").Value = "  -2.76%  "